# Re-style the summary table on slide 16: swap its table style from the
# custom "Table_0" style ({C74D1D18-0C00-411C-9B38-BFF8AB9C3722}) to the
# built-in "Medium Style 2 - Accent 1" style
# ({83BCB8D4-4A02-4C23-8AB7-30548911CDF7}), matching the author's change to
# <a:tableStyleId> in ppt/slides/slide16.xml.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

$table = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $table = $sh.Table
        break
    }
}

if ($table -ne $null) {
    $table.ApplyStyle("{83BCB8D4-4A02-4C23-8AB7-30548911CDF7}")
}
